$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A4: date 2015-12-31 (serial 42369), default (General) number format
$ws.Range("A4").Value = 42369

# New row 6: date 2016-01-06 (serial 42375), times, and a note
# Reuse the same number formats already used elsewhere in the sheet
# (A2 = date format, B2/C2 = time format) so Excel maps to the same style.
$ws.Range("A6").Value = 42375
$ws.Range("A6").NumberFormat = "m/d/yy"

$ws.Range("B6").Value = 0.46875
$ws.Range("B6").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("C6").Value = 0.48958333333333331
$ws.Range("C6").NumberFormat = $ws.Range("C2").NumberFormat

$ws.Range("E6").Value = "Added income tax to payslip."

# Update the active selection to A5
$ws.Range("A5").Select()
